$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 16

$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 14

$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 11

$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 16

$ws.Range("G6").Value = 7
$ws.Range("H6").Value = 33

$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 17

$ws.Range("G8").Value = 5
$ws.Range("H8").Value = 19

$ws.Range("G9").Value = 4
$ws.Range("H9").Value = 10

$ws.Range("E10").Value = 28
$ws.Range("F10").Value = 14
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 16

$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 11

$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 12

$ws.Range("E14").Value = 37
$ws.Range("F14").Value = 17
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 20

$ws.Range("G15").Value = 11
$ws.Range("H15").Value = 55

$ws.Range("E16").Value = 305
$ws.Range("F16").Value = 85
$ws.Range("G16").Value = 87
$ws.Range("H16").Value = 172

$ws.Range("E17").Value = 21
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 12

$ws.Range("G18").Value = 23
$ws.Range("H18").Value = 51
